# closing dates were wrong for algo
# Update intraday HIGH/LOW/LTP/PREV figures on the SBIN sheet to the
# corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SBIN")

# Row 7 (summary row at top: F=open/ref, G=High, H=Low, I=LTP, J=Prev)
$ws.Range("F7").Value = 632.55
$ws.Range("G7").Value = 648.95
$ws.Range("H7").Value = 630.1
$ws.Range("I7").Value = 643.75
$ws.Range("J7").Value = 636.45

# Intraday table rows 9-21 (High Rate / Low Rate / Close Rate columns G/H/I)
$ws.Range("G9").Value = 636
$ws.Range("H9").Value = 627.5
$ws.Range("I9").Value = 632.4

$ws.Range("G10").Value = 638.65
$ws.Range("H10").Value = 631.05
$ws.Range("I10").Value = 634.45

$ws.Range("G11").Value = 639.5
$ws.Range("H11").Value = 633.65
$ws.Range("I11").Value = 636

$ws.Range("G12").Value = 637.9
$ws.Range("H12").Value = 634
$ws.Range("I12").Value = 636.3

$ws.Range("G13").Value = 640
$ws.Range("H13").Value = 634.2
$ws.Range("I13").Value = 634.7

$ws.Range("G14").Value = 638.25
$ws.Range("H14").Value = 634.5
$ws.Range("I14").Value = 637.9

$ws.Range("G15").Value = 640.8
$ws.Range("H15").Value = 637.6
$ws.Range("I15").Value = 640.35

$ws.Range("G16").Value = 643.4
$ws.Range("H16").Value = 638.85
$ws.Range("I16").Value = 639.65

$ws.Range("G17").Value = 640.55
$ws.Range("H17").Value = 638.35
$ws.Range("I17").Value = 639.55

$ws.Range("G18").Value = 640.95
$ws.Range("H18").Value = 637.6
$ws.Range("I18").Value = 638.75

$ws.Range("G19").Value = 640.95
$ws.Range("H19").Value = 636.15
$ws.Range("I19").Value = 640.6

$ws.Range("G20").Value = 648.95
$ws.Range("H20").Value = 640.05
$ws.Range("I20").Value = 642.4

$ws.Range("G21").Value = 645.4
$ws.Range("H21").Value = 642.4
$ws.Range("I21").Value = 644

$wb.Save()
